$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''96.989.56'
$ws.Range('E2').Value = '''  +1.85%  '
$ws.Range('D3').Value = '''3.598.74'
$ws.Range('E3').Value = '''  -0.15%  '
$ws.Range('E4').Value = '''  -0.08%  '
$ws.Range('D5').Value = '''243.36'
$ws.Range('E5').Value = '''  +3.72%  '
$ws.Range('D6').Value = '''657.15'
$ws.Range('E6').Value = '''  +0.60%  '
$ws.Range('E7').Value = '''  +14.17%  '
$ws.Range('E8').Value = '''  +3.15%  '
$ws.Range('E9').Value = '''  +6.60%  '
$ws.Range('E10').Value = '''  -0.04%  '
$ws.Range('D11').Value = '''3.592.39'
$ws.Range('E11').Value = '''  -0.21%  '
$ws.Range('D12').Value = '''43.78'
$ws.Range('E12').Value = '''  +4.41%  '
$ws.Range('E13').Value = '''  +1.76%  '
$ws.Range('D14').Value = '''6.47'
$ws.Range('E14').Value = '''  +1.59%  '
$ws.Range('D15').Value = '''4.263.03'
$ws.Range('E15').Value = '''  -0.82%  '
$ws.Range('D16').Value = '''96.808.86'
$ws.Range('E16').Value = '''  +1.76%  '
$ws.Range('D17').Value = '''0.0000258'
$ws.Range('E17').Value = '''  +2.68%  '
$ws.Range('D18').Value = '''3.583.32'
$ws.Range('E18').Value = '''  -0.45%  '
$ws.Range('B19').Value = '''Polkadot'
$ws.Range('C19').Value = '''https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D19').Value = '''8.09'
$ws.Range('E19').Value = '''  +2.50%  '
$ws.Range('B20').Value = '''Uniswap'
$ws.Range('C20').Value = '''https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').Value = '''12.70'
$ws.Range('E20').Value = '''  -1.55%  '
$ws.Range('D21').Value = '''18.04'
$ws.Range('E21').Value = '''  +1.01%  '
$ws.Range('D22').Value = '''0.531'
$ws.Range('E22').Value = '''  +11.76%  '
$ws.Range('D23').Value = '''511.33'
$ws.Range('E23').Value = '''  +1.10%  '
$ws.Range('D24').Value = '''3.42'
$ws.Range('E24').Value = '''  -2.19%  '
$ws.Range('E25').Value = '''  +3.86%  '
$ws.Range('D26').Value = '''6.87'
$ws.Range('E26').Value = '''  +4.91%  '
$ws.Range('D27').Value = '''97.12'
$ws.Range('E27').Value = '''  +2.21%  '
$ws.Range('D28').Value = '''13.09'
$ws.Range('E28').Value = '''  +5.49%  '
$ws.Range('D29').Value = '''3.785.68'
$ws.Range('E29').Value = '''  -0.42%  '
$ws.Range('D30').Value = '''3.05'
$ws.Range('E30').Value = '''  +0.13%  '
$ws.Range('D31').Value = '''0.150'
$ws.Range('E31').Value = '''  +8.83%  '
$ws.Range('D32').Value = '''11.55'
$ws.Range('E32').Value = '''  +3.66%  '
$ws.Range('D33').Value = '''0.999'
$ws.Range('E33').Value = '''  +0.00%  '
$ws.Range('D34').Value = '''0.186'
$ws.Range('E34').Value = '''  +5.77%  '
$ws.Range('D35').Value = '''0.997'
$ws.Range('E35').Value = '''  +0.00%  '
$ws.Range('D36').Value = '''31.61'
$ws.Range('E36').Value = '''  -2.43%  '
$ws.Range('D37').Value = '''626.91'
$ws.Range('E37').Value = '''  +12.13%  '
$ws.Range('B38').Value = '''PolygonEcosystemToken'
$ws.Range('C38').Value = '''https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D38').Value = '''0.568'
$ws.Range('E38').Value = '''  +2.10%  '
$ws.Range('B39').Value = '''RenderToken'
$ws.Range('C39').Value = '''https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D39').Value = '''8.70'
$ws.Range('E39').Value = '''  +8.23%  '
$ws.Range('D40').Value = '''1.63'
$ws.Range('E40').Value = '''  +11.58%  '
$ws.Range('D41').Value = '''0.154'
$ws.Range('E41').Value = '''  +2.61%  '
$ws.Range('D43').Value = '''1.89'
$ws.Range('E43').Value = '''  +8.49%  '
$ws.Range('D44').Value = '''0.916'
$ws.Range('E44').Value = '''  +0.70%  '
$ws.Range('D45').Value = '''5.89'
$ws.Range('E45').Value = '''  +3.83%  '
$ws.Range('D46').Value = '''0.0430'
$ws.Range('E46').Value = '''  +4.96%  '
$ws.Range('D47').Value = '''2.31'
$ws.Range('E47').Value = '''  +2.57%  '
$ws.Range('D48').Value = '''23.67'
$ws.Range('E48').Value = '''  +0.49%  '
$ws.Range('D49').Value = '''33.24'
$ws.Range('E49').Value = '''  -8.51%  '
$ws.Range('D50').Value = '''8.39'
$ws.Range('E50').Value = '''  +5.23%  '
$ws.Range('E51').Value = '''  -1.76%  '
